$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 21 (id 20): Giuseppe Cangemi, no approvato/numero_tessera/inviato
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Giuseppe"
$ws.Cells.Item(21, 3).Value = "Cangemi"
$ws.Cells.Item(21, 4).Value = "giuseppecangemi94@gmail.com"
$ws.Cells.Item(21, 5).Value = "derryrockfoto.jpg"

# New row 22 (id 21): Pinco Pallino
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "Pinco "
$ws.Cells.Item(22, 3).Value = "Pallino"
$ws.Cells.Item(22, 4).Value = "savvoz@pinko.com"
$ws.Cells.Item(22, 5).Value = "1000052151.jpg"

# New row 23 (id 22): Luca Era Gay
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "Luca "
$ws.Cells.Item(23, 3).Value = "Era Gay"
$ws.Cells.Item(23, 4).Value = "LucaKastlewave@gmail.com"
$ws.Cells.Item(23, 5).Value = "1000052725.jpg"

# New row 24 (id 23): Giuseppe Cangemi, approved with numero_tessera 1209
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "Giuseppe"
$ws.Cells.Item(24, 3).Value = "Cangemi"
$ws.Cells.Item(24, 4).Value = "giuseppecangemi94@gmail.com"
$ws.Cells.Item(24, 5).Value = "derryrockfoto.jpg"
$ws.Cells.Item(24, 6).Value = "SI"
$ws.Cells.Item(24, 7).Value = 1209
$ws.Cells.Item(24, 8).Value = "SI"
